$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.367.03'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.395.23'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.46%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '549.27'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.19'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.537'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -10.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.394.48'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.30'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.348'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.43'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.824.57'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000166'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.050.06'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.393.15'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.81'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.15'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '318.69'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.73'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.97%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.92'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.76'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +8.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '540.11'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.53%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0943'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.510.81'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.45'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.10'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.146'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.44%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.58'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.56'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.73'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.377'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.86'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.13'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '138.44'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -7.18%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.25'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.21'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '141.72'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.63'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.30'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0521'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.579'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.09%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.52%  '
